# Update cryptos list with refreshed prices/volumes and reshuffled rankings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.750.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.340.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.36"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.339.17"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.763.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.634.50"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.343.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.69"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.50"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.03%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.456.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.35%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "496.51"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.37"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.33%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.67%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.146"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.97%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.60"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.63%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.24"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.69%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.70%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "142.13"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.44%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.00%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.93%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.47%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.37%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.567"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0900"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.57%  "

Write-Output "Applied cryptos update."
